$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'62.851.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.353.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'572.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'152.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.56%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "'3.355.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.84%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'7.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.118"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.51%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.96%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'3.930.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.89%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.03%  "

# Row 15 - was -> ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.69%  "

# Row 16 - was -> Avalanche
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'26.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'62.869.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.11%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.360.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.84%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'13.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.93%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'8.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'385.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.93%  "

# Row 23 - Dai
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "

# Row 24 - Polygon
$ws.Range("D24").Value = "'0.536"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.93%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'70.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'9.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.47%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +6.65%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "'0.0₃0971"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.26%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.35%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.99%  "

# Row 31 - was -> Fetch.AI
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.52%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  +5.69%  "

# Row 33 - was -> RenderToken
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "'6.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.83%  "

# Row 34 - was -> EthereumClassic
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.84%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'6.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.48%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.33%  "

# Row 37 - Monero
$ws.Range("D37").Value = "'158.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +12.61%  "

# Row 39 - EnergySwap
$ws.Range("D39").Value = "'27.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.28%  "

# Row 40 - was -> VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0330"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.01%  "

# Row 41 - was -> Hedera
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0739"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.80%  "

# Row 42 - Maker
$ws.Range("D42").Value = "'2.785.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'41.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.60%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'4.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.02%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "'0.745"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.73%  "

# Row 46 - ONDO
$ws.Range("D46").Value = "'1.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.90%  "

# Row 47 - was -> RenzoRestakedETH
$ws.Range("B47").Value = "RenzoRestakedETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D47").Value = "'3.397.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.04%  "

# Row 48 - was -> InjectiveProtocol
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'22.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.05%  "

# Row 49 - was -> Cosmos
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.80%  "

# Row 50 - was -> Stellar
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51 - Bittensor
$ws.Range("D51").Value = "'290.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.50%  "
